# live_trading_results.xlsx - Trade #34 closed at 2026-02-17 20:57:51 - unknown UNKNOWN +0.000%
#
# This applies:
#  1) Updated rollup numbers on "Summary"
#  2) Updated MarketMaking rollup row on "Strategy Status"
#  3) Closes the previously-OPEN MarketMaking trade (early_exit) on
#     "All Trades" (row 63) and on "MarketMaking" (row 30)
#  4) Appends the freshly-opened MarketMaking trade as a new row on both
#     "All Trades" (row 96) and "MarketMaking" (row 63)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(3, 2).Value = 1400.45   # Current Capital
$summary.Cells.Item(4, 2).Value = 0.24      # Total P&L $
$summary.Cells.Item(6, 2).Value = 62        # Total Trades
$summary.Cells.Item(8, 2).Value = 25        # Losing Trades
$summary.Cells.Item(9, 2).Value = 46.77     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Cells.Item(5, 3).Value = 100.45   # Capital
$status.Cells.Item(5, 4).Value = 29       # Trades
$status.Cells.Item(5, 5).Value = 0.13     # P&L $
$status.Cells.Item(5, 6).Value = 0.45     # P&L %
$status.Cells.Item(5, 7).Value = 51.72    # Win Rate %

# ---------------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close the existing open trade (row 63, Trade # 62) as an early_exit
$allTrades.Cells.Item(63, 7).Value = 0.08        # Exit Price
$allTrades.Cells.Item(63, 8).Value = "CLOSED"    # Status
$allTrades.Cells.Item(63, 9).Value = -11.1111    # P&L %
$allTrades.Cells.Item(63, 10).Value = -0.01      # P&L $
$allTrades.Cells.Item(63, 11).Value = 100.45     # Capital After
$allTrades.Cells.Item(63, 12).Value = "early_exit" # Exit Reason
$allTrades.Cells.Item(63, 13).Value = 0.11       # Duration (min)

# Append the newly-opened trade as row 96 (Trade # 95)
$allTrades.Cells.Item(96, 1).Value = 95
# Copy the date cell from the row above so the "2026-02-17" text isn't
# auto-converted into a real date value/format by Excel's type inference.
$allTrades.Cells.Item(95, 2).Copy($allTrades.Cells.Item(96, 2))
$allTrades.Cells.Item(96, 3).Value = "20:57:44"
$allTrades.Cells.Item(96, 4).Value = "MarketMaking"
$allTrades.Cells.Item(96, 5).Value = "UP"
$allTrades.Cells.Item(96, 6).Value = 0.09
$allTrades.Cells.Item(96, 8).Value = "OPEN"
$allTrades.Cells.Item(96, 9).Value = 0
$allTrades.Cells.Item(96, 10).Value = 0
$allTrades.Cells.Item(96, 11).Value = 100.4555022889912
$allTrades.Cells.Item(96, 13).Value = 0
$allTrades.Cells.Item(96, 14).Value = 0
$allTrades.Cells.Item(96, 15).Value = 0
$allTrades.Cells.Item(96, 16).Value = 0.6
$allTrades.Cells.Item(96, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# 4) MarketMaking sheet
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Close the existing open trade (row 30, Trade # 62) as an early_exit
$mm.Cells.Item(30, 7).Value = 0.08          # Exit Price
$mm.Cells.Item(30, 8).Value = "CLOSED"      # Status
$mm.Cells.Item(30, 9).Value = -11.1111      # P&L %
$mm.Cells.Item(30, 10).Value = -0.01        # P&L $
$mm.Cells.Item(30, 11).Value = 100.45       # Capital After
$mm.Cells.Item(30, 16).Value = "early_exit" # Exit Reason
$mm.Cells.Item(30, 17).Value = 0.11         # Duration (min)

# Append the newly-opened trade as row 63 (Trade # 95)
$mm.Cells.Item(63, 1).Value = 95
# Copy the date cell from the row above so the "2026-02-17" text isn't
# auto-converted into a real date value/format by Excel's type inference.
$mm.Cells.Item(62, 2).Copy($mm.Cells.Item(63, 2))
$mm.Cells.Item(63, 3).Value = "20:57:44"
$mm.Cells.Item(63, 4).Value = "MarketMaking"
$mm.Cells.Item(63, 5).Value = "UP"
$mm.Cells.Item(63, 6).Value = 0.09
$mm.Cells.Item(63, 8).Value = "OPEN"
$mm.Cells.Item(63, 9).Value = 0
$mm.Cells.Item(63, 10).Value = 0
$mm.Cells.Item(63, 11).Value = 100.4555022889912
$mm.Cells.Item(63, 12).Value = 0
$mm.Cells.Item(63, 13).Value = 0
$mm.Cells.Item(63, 14).Value = 0.6
$mm.Cells.Item(63, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(63, 17).Value = 0
